$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5700.5
$ws.Range("J86").Value = 8261.1
$ws.Range("L86").Value = 8261.1
$ws.Range("N86").Value = -10507.1
$ws.Range("H89").Value = 5700.5
$ws.Range("J89").Value = 8261.1
$ws.Range("L89").Value = 41305.5
$ws.Range("N89").Value = -52537.5
$ws.Range("H113").Value = 5214.2856
$ws.Range("J113").Value = 5760
$ws.Range("L113").Value = 5760
$ws.Range("N113").Value = -12268
$ws.Range("H127").Value = 6265
$ws.Range("I127").Value = 6265
$ws.Range("K127").Value = 18795
$ws.Range("M127").Value = -13835
$ws.Range("H132").Value = 3625.9333
$ws.Range("I132").Value = 3621.4075
$ws.Range("K132").Value = 10864.2225
$ws.Range("M132").Value = -8334.2225
$ws.Range("H138").Value = 3807.8718
$ws.Range("J138").Value = 6750
$ws.Range("L138").Value = 20250
$ws.Range("N138").Value = -30530
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2627
$ws.Range("I2").Value = 1502.6666
$ws.Range("K2").Value = 1502.6666
$ws.Range("M2").Value = -1389.6666
$ws.Range("H5").Value = 269.16666
$ws.Range("J5").Value = 71.5
$ws.Range("L5").Value = 71.5
$ws.Range("N5").Value = -295.5
$ws.Range("H32").Value = 5081.2603
$ws.Range("I32").Value = 1984.6029
$ws.Range("K32").Value = 1984.6029
$ws.Range("M32").Value = -1697.6029
$ws.Range("H45").Value = 1483.4865
$ws.Range("I45").Value = 1290.9706
$ws.Range("J45").Value = 3665.3333
$ws.Range("K45").Value = 1290.9706
$ws.Range("L45").Value = 3665.3333
$ws.Range("M45").Value = -913.9706000000001
$ws.Range("N45").Value = -4419.3333
$ws.Range("H116").Value = 2627
$ws.Range("I116").Value = 1502.6666
$ws.Range("K116").Value = 1502.6666
$ws.Range("M116").Value = 791.3334
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2627
$ws.Range("I3").Value = 1502.6666
$ws.Range("K3").Value = 1502.6666
$ws.Range("M3").Value = -1388.6666
$ws.Range("H4").Value = 269.16666
$ws.Range("J4").Value = 71.5
$ws.Range("L4").Value = 71.5
$ws.Range("N4").Value = -301.5
$ws.Range("H86").Value = 2064.2173
$ws.Range("I86").Value = 1774.4
$ws.Range("K86").Value = 1774.4
$ws.Range("M86").Value = -651.4000000000001
$ws.Range("H89").Value = 2064.2173
$ws.Range("I89").Value = 1774.4
$ws.Range("K89").Value = 8872
$ws.Range("M89").Value = -3256
$ws.Range("H105").Value = 1487.04
$ws.Range("I105").Value = 1410.125
$ws.Range("J105").Value = 3333
$ws.Range("K105").Value = 1410.125
$ws.Range("L105").Value = 3333
$ws.Range("M105").Value = 336.875
$ws.Range("N105").Value = -6827
$ws.Range("H107").Value = 3333
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 3333
$ws.Range("K107").Value = 0
$ws.Range("M107").Value = 3333
$ws.Range("N107").Value = -7173
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9310.5
$ws.Range("I31").Value = 4373.3335
$ws.Range("K31").Value = 4373.3335
$ws.Range("M31").Value = -4078.3335
$ws.Range("H34").Value = 9310.5
$ws.Range("I34").Value = 4373.3335
$ws.Range("K34").Value = 4373.3335
$ws.Range("M34").Value = -4171.3335
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H53").Value = 70000
$ws.Range("J53").Value = 70000
$ws.Range("L53").Value = 70000
$ws.Range("N53").Value = -71214
$ws.Range("H99").Value = 2278.1904
$ws.Range("I99").Value = 1883.125
$ws.Range("K99").Value = 1883.125
$ws.Range("M99").Value = -385.125
$ws.Range("H105").Value = 1728.5625
$ws.Range("I105").Value = 1900.1538
$ws.Range("J105").Value = 985
$ws.Range("K105").Value = 1900.1538
$ws.Range("L105").Value = 985
$ws.Range("M105").Value = -153.1538
$ws.Range("N105").Value = -4479
$ws.Range("H106").Value = 255357.14
$ws.Range("J106").Value = 255357.14
$ws.Range("L106").Value = 255357.14
$ws.Range("N106").Value = -257881.14
$ws.Range("H107").Value = 1996.2142
$ws.Range("I107").Value = 1837.4
$ws.Range("J107").Value = 2393.25
$ws.Range("K107").Value = 1837.4
$ws.Range("L107").Value = 2393.25
$ws.Range("M107").Value = 82.59999999999991
$ws.Range("N107").Value = -6233.25
$ws.Range("H126").Value = 2278.1904
$ws.Range("I126").Value = 1883.125
$ws.Range("K126").Value = 5649.375
$ws.Range("M126").Value = -3179.375
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2900.5
$ws.Range("I5").Value = 2900.5
$ws.Range("K5").Value = 8701.5
$ws.Range("M5").Value = -8589.5
$ws.Range("H56").Value = 4013.8386
$ws.Range("I56").Value = 4013.8386
$ws.Range("K56").Value = 4013.8386
$ws.Range("M56").Value = -3483.8386
$ws.Range("H61").Value = 148.75
$ws.Range("I61").Value = 148.75
$ws.Range("K61").Value = 446.25
$ws.Range("M61").Value = -231.25
$ws.Range("H135").Value = 2900.5
$ws.Range("I135").Value = 2900.5
$ws.Range("K135").Value = 26104.5
$ws.Range("M135").Value = -23569.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 5333.3335
$ws.Range("I31").Value = 5400
$ws.Range("K31").Value = 5400
$ws.Range("M31").Value = -5108
$ws.Range("H37").Value = 5333.3335
$ws.Range("I37").Value = 5400
$ws.Range("K37").Value = 5400
$ws.Range("M37").Value = -5123
$ws.Range("H41").Value = 5000
$ws.Range("I41").Value = 5000
$ws.Range("K41").Value = 5000
$ws.Range("M41").Value = -4645
$ws.Range("H70").Value = 4233.2
$ws.Range("I70").Value = 2246.5
$ws.Range("K70").Value = 2246.5
$ws.Range("M70").Value = -1976.5
$ws.Range("H73").Value = 4233.2
$ws.Range("I73").Value = 2246.5
$ws.Range("K73").Value = 2246.5
$ws.Range("M73").Value = -1310.5
$ws.Range("H113").Value = 4217.7
$ws.Range("J113").Value = 5072.8335
$ws.Range("L113").Value = 5072.8335
$ws.Range("N113").Value = -9412.833500000001
$ws.Range("H126").Value = 4284.143
$ws.Range("I126").Value = 1997.25
$ws.Range("J126").Value = 7333.3335
$ws.Range("K126").Value = 5991.75
$ws.Range("L126").Value = 22000.0005
$ws.Range("M126").Value = -3521.75
$ws.Range("N126").Value = -26940.0005
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3522.1904
$ws.Range("I40").Value = 2344.2
$ws.Range("K40").Value = 2344.2
$ws.Range("M40").Value = -2208.2
$ws.Range("H122").Value = 5501.6313
$ws.Range("I122").Value = 6352.077
$ws.Range("J122").Value = 3659
$ws.Range("K122").Value = 19056.231
$ws.Range("L122").Value = 10977
$ws.Range("M122").Value = -16606.231
$ws.Range("N122").Value = -15877
$ws.Range("H132").Value = 4670.5405
$ws.Range("I132").Value = 4285.2256
$ws.Range("K132").Value = 12855.6768
$ws.Range("M132").Value = -10325.6768
$ws.Range("H136").Value = 4696.7837
$ws.Range("I136").Value = 3581.8262
$ws.Range("K136").Value = 10745.4786
$ws.Range("M136").Value = -8195.4786
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1540438.9
$ws.Range("I107").Value = 2502071
$ws.Range("J107").Value = 1827.4
$ws.Range("K107").Value = 7506213
$ws.Range("L107").Value = 5482.200000000001
$ws.Range("M107").Value = -7504293
$ws.Range("N107").Value = -9322.200000000001
$ws.Range("H119").Value = 178812.5
$ws.Range("J119").Value = 178812.5
$ws.Range("L119").Value = 178812.5
$ws.Range("N119").Value = -188488.5
$ws.Range("H122").Value = 5757.4443
$ws.Range("I122").Value = 4474.6665
$ws.Range("J122").Value = 8323
$ws.Range("K122").Value = 13423.9995
$ws.Range("L122").Value = 24969
$ws.Range("M122").Value = -10973.9995
$ws.Range("N122").Value = -29869
$ws.Range("H128").Value = 299950
$ws.Range("J128").Value = 299950
$ws.Range("L128").Value = 299950
$ws.Range("N128").Value = -309910
$ws.Range("H132").Value = 5377.875
$ws.Range("I132").Value = 1927.4615
$ws.Range("J132").Value = 20329.666
$ws.Range("K132").Value = 5782.3845
$ws.Range("L132").Value = 60988.99800000001
$ws.Range("M132").Value = -3252.3845
$ws.Range("N132").Value = -66048.99800000001
$ws.Range("H136").Value = 3966.889
$ws.Range("I136").Value = 2865.9546
$ws.Range("K136").Value = 8597.863799999999
$ws.Range("M136").Value = -6047.863799999999
